# Weekly update: insert the newest "Coliflor" market entries (Primera /
# Segunda) for Femacal de La Calera, pushing the existing history down by
# two rows (the sheet is ordered most-recent-first).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right below the current top data rows (491:492),
# shifting all existing data (old rows 491:521) down to 493:523.
$ws.Rows("491:492").Insert()

# New row 491: Primera quality, fecha 2022-01-24 (serial 44585)
$ws.Cells.Item(491, 1).Value = 3
$ws.Cells.Item(491, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(491, 3).Value = "Coquimbo"
$ws.Cells.Item(491, 4).Value = 44585
$ws.Cells.Item(491, 5).Value = 5
$ws.Cells.Item(491, 6).Value = 100112008
$ws.Cells.Item(491, 7).Value = "Coliflor"
$ws.Cells.Item(491, 8).Value = "Sin especificar"
$ws.Cells.Item(491, 9).Value = "Primera"
$ws.Cells.Item(491, 10).Value = 850
$ws.Cells.Item(491, 11).Value = 1300
$ws.Cells.Item(491, 12).Value = 1300
$ws.Cells.Item(491, 13).Value = 1300
$ws.Cells.Item(491, 14).Value = "$/unidad"
$ws.Cells.Item(491, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(491, 16).Value = 1300
$ws.Cells.Item(491, 17).Value = 1
$ws.Cells.Item(491, 18).Value = "Hortaliza"

# New row 492: Segunda quality, same fecha
$ws.Cells.Item(492, 1).Value = 3
$ws.Cells.Item(492, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(492, 3).Value = "Coquimbo"
$ws.Cells.Item(492, 4).Value = 44585
$ws.Cells.Item(492, 5).Value = 5
$ws.Cells.Item(492, 6).Value = 100112008
$ws.Cells.Item(492, 7).Value = "Coliflor"
$ws.Cells.Item(492, 8).Value = "Sin especificar"
$ws.Cells.Item(492, 9).Value = "Segunda"
$ws.Cells.Item(492, 10).Value = 950
$ws.Cells.Item(492, 11).Value = 1000
$ws.Cells.Item(492, 12).Value = 1000
$ws.Cells.Item(492, 13).Value = 1000
$ws.Cells.Item(492, 14).Value = "$/unidad"
$ws.Cells.Item(492, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(492, 16).Value = 1000
$ws.Cells.Item(492, 17).Value = 1
$ws.Cells.Item(492, 18).Value = "Hortaliza"
